$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 407; this shifts the existing rows 407-434
# down to 408-435 (matching the diff's row-shift pattern) and bumps the
# sheet dimension to A1:R435 automatically.
$ws.Rows.Item(407).Insert()

# Populate the newly inserted row 407 with the new weekly record.
$ws.Range("A407").Value = 4
$ws.Range("B407").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C407").Value = "Los Lagos"
$ws.Range("D407").Value = 45265
$ws.Range("D407").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E407").Value = 10
$ws.Range("F407").Value = 100112039
$ws.Range("G407").Value = "Ciboulette"
$ws.Range("H407").Value = "Sin especificar"
$ws.Range("I407").Value = "Primera"
$ws.Range("J407").Value = 240
$ws.Range("K407").Value = 2500
$ws.Range("L407").Value = 2500
$ws.Range("M407").Value = 2500
$ws.Range("N407").Value = '$/docena de atados'
$ws.Range("O407").Value = "Región Metropolitana"
$ws.Range("P407").Value = 833
$ws.Range("Q407").Value = 3
$ws.Range("R407").Value = "Hortaliza"
